$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
# F8: was green-filled thick-border (style 5) value 2 -> plain thick-border (style 2) value 5
$ws.Range("G4").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = 5

# G8, H8: already style 2 (plain thick border), just fill in the grade
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 5

# I8, J8: new cells, style 7 (no fill, thick left/right border)
$ws.Range("I6").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = 5

$ws.Range("J6").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J8").Value = 5

# --- Row 12 ---
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 5

# I12, J12: new cells, style 6 (green fill, thick left/right border)
$ws.Range("I13").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 5

$ws.Range("J13").Copy()
$ws.Range("J12").PasteSpecial(-4122)
$ws.Range("J12").Value = 5

# --- Row 14 ---
# I14, J14: new cells, no explicit style
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 5

# --- Row 20 ---
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 5

$ws.Range("I13").Copy()
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("I20").Value = 5

$ws.Range("J13").Copy()
$ws.Range("J20").PasteSpecial(-4122)
$ws.Range("J20").Value = 5

# --- Row 23 ---
$ws.Range("H23").Value = 5

$ws.Range("I13").Copy()
$ws.Range("I23").PasteSpecial(-4122)
$ws.Range("I23").Value = 5

$ws.Range("J13").Copy()
$ws.Range("J23").PasteSpecial(-4122)
$ws.Range("J23").Value = 5

$excel.CutCopyMode = $false

# --- View state: frozen pane scrolled back up, selection moved to I12 ---
$ws.Range("I12").Select() | Out-Null
